$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-05 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-06 Monday", 2)
$d.Content.Find.Execute("535÷3=178, 1", $true, $false, $false, $false, $false, $true, 1, $false, "517÷7=73, 6", 2)
$d.Content.Find.Execute("364÷2=182, 0", $true, $false, $false, $false, $false, $true, 1, $false, "844÷7=120, 4", 2)
$d.Content.Find.Execute("564÷9=62, 6", $true, $false, $false, $false, $false, $true, 1, $false, "860÷5=172, 0", 2)
$d.Content.Find.Execute("731÷3=243, 2", $true, $false, $false, $false, $false, $true, 1, $false, "615÷5=123, 0", 2)
$d.Content.Find.Execute("289÷2=144, 1", $true, $false, $false, $false, $false, $true, 1, $false, "927÷8=115, 7", 2)
$d.Content.Find.Execute("936÷5=187, 1", $true, $false, $false, $false, $false, $true, 1, $false, "649÷2=324, 1", 2)
$d.Content.Find.Execute("973÷3=324, 1", $true, $false, $false, $false, $false, $true, 1, $false, "153÷2=76, 1", 2)
$d.Content.Find.Execute("485÷2=242, 1", $true, $false, $false, $false, $false, $true, 1, $false, "331÷7=47, 2", 2)
$d.Content.Find.Execute("840÷9=93, 3", $true, $false, $false, $false, $false, $true, 1, $false, "446÷2=223, 0", 2)
$d.Content.Find.Execute("427÷8=53, 3", $true, $false, $false, $false, $false, $true, 1, $false, "464÷8=58, 0", 2)
$d.Content.Find.Execute("282÷6=47, 0", $true, $false, $false, $false, $false, $true, 1, $false, "951÷3=317, 0", 2)
$d.Content.Find.Execute("135÷3=45, 0", $true, $false, $false, $false, $false, $true, 1, $false, "745÷8=93, 1", 2)
$d.Content.Find.Execute("190÷8=23, 6", $true, $false, $false, $false, $false, $true, 1, $false, "207÷2=103, 1", 2)
$d.Content.Find.Execute("852÷8=106, 4", $true, $false, $false, $false, $false, $true, 1, $false, "445÷5=89, 0", 2)
$d.Content.Find.Execute("440÷3=146, 2", $true, $false, $false, $false, $false, $true, 1, $false, "126÷7=18, 0", 2)
$d.Content.Find.Execute("813÷9=90, 3", $true, $false, $false, $false, $false, $true, 1, $false, "654÷8=81, 6", 2)
$d.Content.Find.Execute("425÷4=106, 1", $true, $false, $false, $false, $false, $true, 1, $false, "762÷6=127, 0", 2)
$d.Content.Find.Execute("803÷9=89, 2", $true, $false, $false, $false, $false, $true, 1, $false, "393÷2=196, 1", 2)
$d.Content.Find.Execute("842÷2=421, 0", $true, $false, $false, $false, $false, $true, 1, $false, "394÷6=65, 4", 2)
$d.Content.Find.Execute("586÷7=83, 5", $true, $false, $false, $false, $false, $true, 1, $false, "510÷7=72, 6", 2)
$d.Content.Find.Execute("647÷5=129, 2", $true, $false, $false, $false, $false, $true, 1, $false, "688÷2=344, 0", 2)
$d.Content.Find.Execute("482÷8=60, 2", $true, $false, $false, $false, $false, $true, 1, $false, "782÷4=195, 2", 2)
$d.Content.Find.Execute("747÷5=149, 2", $true, $false, $false, $false, $false, $true, 1, $false, "585÷8=73, 1", 2)
$d.Content.Find.Execute("648÷7=92, 4", $true, $false, $false, $false, $false, $true, 1, $false, "199÷9=22, 1", 2)
$d.Content.Find.Execute("310÷5=62, 0", $true, $false, $false, $false, $false, $true, 1, $false, "407÷7=58, 1", 2)
